# Update "想去人数" (number of people interested) figures on the
# "展览" and "全部类型" worksheets to reflect the freshly generated data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 5015
    $ws.Range("F3").Value = 152
    $ws.Range("F4").Value = 883
}
